# Update "想去人数" (F column) values on sheets "展览" and "全部类型"
# to reflect refreshed counts from the latest data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 561
$ws1.Range("F4").Value  = 206
$ws1.Range("F6").Value  = 516
$ws1.Range("F7").Value  = 109
$ws1.Range("F9").Value  = 50
$ws1.Range("F10").Value = 6871
$ws1.Range("F11").Value = 241
$ws1.Range("F12").Value = 382
$ws1.Range("F13").Value = 3222
$ws1.Range("F14").Value = 219
$ws1.Range("F15").Value = 388
$ws1.Range("F17").Value = 559
$ws1.Range("F18").Value = 30

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 561
$ws4.Range("F6").Value  = 206
$ws4.Range("F8").Value  = 516
$ws4.Range("F9").Value  = 109
$ws4.Range("F11").Value = 50
$ws4.Range("F13").Value = 6871
$ws4.Range("F15").Value = 241
$ws4.Range("F16").Value = 382
$ws4.Range("F17").Value = 3222
$ws4.Range("F18").Value = 219
$ws4.Range("F19").Value = 388
$ws4.Range("F21").Value = 559
$ws4.Range("F22").Value = 30
